# Appends rows 179-186 to the "PODSUMOWANIE" sheet (OLX monitor run 2026-02-22 20:47),
# mirroring the row-178 layout/styling for the new entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

$ws.Range("A179").Value2 = "2026-02-22 20:47:59"
$ws.Range("B179").Value2 = "poqui"
$ws.Range("C179").Value2 = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza"
$ws.Range("D179").Value2 = 2049
$ws.Range("E179").Value2 = "19.12.2025"
$ws.Range("F179").Value2 = 65
$ws.Range("G179").Value2 = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html"
$ws.Range("H179").Value2 = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"

$ws.Range("A180").Value2 = "2026-02-22 20:47:59"
$ws.Range("B180").Value2 = "poqui"
$ws.Range("C180").Value2 = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda"
$ws.Range("D180").Value2 = 2299
$ws.Range("E180").Value2 = "19.01.2026"
$ws.Range("F180").Value2 = 34
$ws.Range("G180").Value2 = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html"
$ws.Range("H180").Value2 = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"

$ws.Range("A181").Value2 = "2026-02-22 20:47:59"
$ws.Range("B181").Value2 = "poqui"
$ws.Range("C181").Value2 = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"
$ws.Range("D181").Value2 = 2499
$ws.Range("E181").Value2 = "28.10.2025"
$ws.Range("F181").Value2 = 117
$ws.Range("G181").Value2 = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"
$ws.Range("H181").Value2 = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"

$ws.Range("A182").Value2 = "2026-02-22 20:47:59"
$ws.Range("B182").Value2 = "poqui"
$ws.Range("C182").Value2 = "Przytulny pokój blisko Politechniki – ul. Przytulna"
$ws.Range("D182").Value2 = 549
$ws.Range("E182").Value2 = "'10.10.2025"
$ws.Range("F182").Value2 = 135
$ws.Range("G182").Value2 = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"
$ws.Range("H182").Value2 = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"

$ws.Range("A183").Value2 = "2026-02-22 20:47:59"
$ws.Range("B183").Value2 = "pokojewlublinie"
$ws.Range("C183").Value2 = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"
$ws.Range("D183").Value2 = 0
$ws.Range("E183").Value2 = "'11.08.2025"
$ws.Range("F183").Value2 = 195
$ws.Range("G183").Value2 = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"
$ws.Range("H183").Value2 = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"

$ws.Range("A184").Value2 = "2026-02-22 20:47:59"
$ws.Range("B184").Value2 = "pokojewlublinie"
$ws.Range("C184").Value2 = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"
$ws.Range("D184").Value2 = 12640
$ws.Range("E184").Value2 = "19.01.2026"
$ws.Range("F184").Value2 = 34
$ws.Range("G184").Value2 = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"
$ws.Range("H184").Value2 = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"

$ws.Range("A185").Value2 = "2026-02-22 20:47:59"
$ws.Range("B185").Value2 = "dawnypatron"
$ws.Range("C185").Value2 = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."
$ws.Range("D185").Value2 = 730
$ws.Range("E185").Value2 = "20.09.2024"
$ws.Range("F185").Value2 = 520
$ws.Range("G185").Value2 = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"
$ws.Range("H185").Value2 = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"

$ws.Range("A186").Value2 = "2026-02-22 20:47:59"
$ws.Range("B186").Value2 = "dawnypatron"
$ws.Range("C186").Value2 = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"
$ws.Range("D186").Value2 = 14690
$ws.Range("E186").Value2 = "'05.12.2025"
$ws.Range("F186").Value2 = 79
$ws.Range("G186").Value2 = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"
$ws.Range("H186").Value2 = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"

# Re-stamp the new rows with row 178's cell styles (tiles A178:H178 across A179:H186)
$ws.Range("A178:H178").Copy()
$ws.Range("A179:H186").PasteSpecial(-4122)

# Rows 180 and 184 use the plain (non-highlighted) F-column style, like F7/F10
$ws.Range("F7").Copy()
$ws.Range("F180").PasteSpecial(-4122)
$ws.Range("F184").PasteSpecial(-4122)

$excel.CutCopyMode = 0
Write-Host "Added rows 179-186 to PODSUMOWANIE"
